# 40th commit Final 001-048
# - add a new test-case sheet "tc_047" (after "tc_011") holding the
#   reject-comment-as-V4-user automation test data
# - that new sheet becomes the active/selected sheet
# - the previously-active sheet ("Basic_Information_Release") is no longer
#   the active tab (handled automatically once a different sheet is active)
# - the selection on "tc_001" moves from C17 to D1:D2

$wb = $excel.ActiveWorkbook

# --- update the remembered selection on tc_001 --------------------------
$tc001 = $wb.Worksheets.Item("tc_001")
$tc001.Range("D1:D2").Select() | Out-Null

# --- add the new tc_047 sheet, after the last existing sheet (tc_011) ---
$template = $wb.Worksheets.Item("tc_011")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy([System.Reflection.Missing]::Value, $lastSheet) | Out-Null

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "tc_047"

# tc_011's layout has two columns (A/B); tc_047 only needs column A
$newSheet.Range("B1:B2").Clear() | Out-Null

$newSheet.Range("A1").Value = "reject_comment_v4_user"
$newSheet.Range("A2").Value = "Reject comment as V4 user automation test"

$newSheet.Columns.Item(1).ColumnWidth = 41

# leave the new sheet active, with A2 selected
$newSheet.Range("A2").Select() | Out-Null
